# Update the "Configuration" sheet's mapping/path templates so that:
#  - Dataset path template changes from "{YYYY}{MM}{DD}" to "**"
#  - Next path template changes from "{Label}" to "{Directory}"
#  - Session mapping recipe changes from "YYYY"/"MM"/"DD" (three cells) to a
#    single "DICOM:StudyDate" cell (D8/E8 are cleared)
#  - Dataset mapping recipe changes from "Label" to "Directory"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# Row 2: Paths
$ws.Range("D2").Value = "**"
$ws.Range("E2").Value = "{Directory}"

# Row 8: Session mapping - replace YYYY/MM/DD split with single DICOM:StudyDate
$ws.Range("C8").Value = "DICOM:StudyDate"
$ws.Range("D8").Value = $null
$ws.Range("E8").Value = $null

# Row 9: Dataset mapping - Label -> Directory
$ws.Range("C9").Value = "Directory"
